# Bond dates update: recompute "days since previous payment" (col G) and
# "days until next payment" (col I) as if one day has elapsed, while the
# underlying payment dates themselves (col F / col H) stay the same.
#   - Dni od poprzedniej wypłaty (G): +1 day (only where a value exists)
#   - Dni do następnej wypłaty   (I): -1 day (present on every data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count   # header on row 1, data rows 2..lastRow

for ($r = 2; $r -le $lastRow; $r++) {
    $gCell = $ws.Cells.Item($r, 7)   # column G
    $gVal = $gCell.Value()
    if ($gVal -ne $null) {
        $gCell.Value = $gVal + 1
    }

    $iCell = $ws.Cells.Item($r, 9)   # column I
    $iVal = $iCell.Value()
    if ($iVal -ne $null) {
        $iCell.Value = $iVal - 1
    }
}
